$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 59; this shifts existing rows 59-75 down to 60-76
$ws.Rows.Item(59).Insert()

# Populate the new row 59 with the new weekly record
$ws.Cells.Item(59, 1).Value = 4
$ws.Cells.Item(59, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value = "Los Lagos"
$ws.Cells.Item(59, 4).Value = 45218
$ws.Cells.Item(59, 5).Value = 10
$ws.Cells.Item(59, 6).Value = 300000000
$ws.Cells.Item(59, 7).Value = "Espárragos"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 200
$ws.Cells.Item(59, 11).Value = 1900
$ws.Cells.Item(59, 12).Value = 2000
$ws.Cells.Item(59, 13).Value = 1950
$ws.Cells.Item(59, 14).Value = "`$/kilo"
$ws.Cells.Item(59, 15).Value = "Provincia de Linares"
$ws.Cells.Item(59, 16).Value = 1950
$ws.Cells.Item(59, 17).Value = 1
$ws.Cells.Item(59, 18).Value = "Hortaliza"
